# Inventory Edit functionaly finished
# Append 3 new contract-WIP rows (103-105) to the bottom of the tracking
# sheet, mirroring the existing row layout (10 text columns: A..J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then the 10 column values (A..J) as strings.
$newRows = @(
    @{
        Row = 103
        Values = @(
            "4/24/2019",
            "SPE7L4-19-V-0874",
            "42",
            "`$2,389.38 ",
            "5935012860663",
            "CONNECTOR,PLUG,ELEC",
            "ITT",
            "CIR06G2-14S-6S-F80",
            "ZZ",
            "2019 OCT 11"
        )
    },
    @{
        Row = 104
        Values = @(
            "4/24/2019",
            "SPE8E8-19-V-1727",
            "6",
            "`$5,114.82",
            "4420015045491",
            "BAFFLE,FLUID COOLER",
            "Fluid Handling",
            "4371641-0476",
            "CP",
            "2019 OCT 01"
        )
    },
    @{
        Row = 105
        Values = @(
            "4/24/2019",
            "SPE7M1-19-V-6499",
            "23",
            "`$1,677.39 ",
            "5999012232712",
            "DELAY LINE",
            "Data Delay Devices",
            "DDU7-8212",
            "CP",
            "2019 OCT 01"
        )
    }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($r, $col)
        $val = $vals[$i]

        # Columns hold plain text in the source data (dates, quantities,
        # dollar totals and NSNs are all stored as text, not real
        # numbers/dates). Values that *look* numeric/date-like would
        # otherwise get auto-converted by Excel's type inference, so force
        # those specific cells to Text before writing them. Values that are
        # already unambiguous (ids, names, part numbers, etc.) are written
        # as-is and naturally remain text.
        $looksNumeric = $val -match '^\s*-?\$?[0-9][0-9,]*(\.[0-9]+)?\s*$' -or `
                         $val -match '^\d{1,2}/\d{1,2}/\d{2,4}$'

        if ($looksNumeric) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}

Write-Output "Appended rows 103-105 to $($ws.Name)"
